# Fix Training Data Issue (#48)
#
# The "Date" column (BF) stored the value used to derive the source
# filename ("4-29-2013-14" -> month-day + season label) instead of the
# actual calendar date the games were played on. NBA.com reports late
# April games under the "2013-14" season label even though the games
# themselves were played in 2014, which threw off downstream model
# training. Replace the mislabeled text with the correct ISO date
# (2014-04-29) for every data row, leaving the header row and all other
# cells/formatting untouched.
#
# Assigning a "YYYY-MM-DD"-looking string straight to Range.Value makes
# Excel "helpfully" parse it as a date and reformat the cell, which the
# source data does not want (the column holds plain text). To store the
# literal text without touching the cell's style, stage the text (forced
# to the Text format) in an unused scratch cell, copy it, and paste only
# the *value* into the destination - then restore the scratch cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "4-29-2013-14"
$newValue = "2014-04-29"
$dateCol = 58   # column BF

$scratch = $ws.Range("A1")   # unused cell, outside the populated columns

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $dst = $ws.Cells.Item($r, $dateCol)
    if ($dst.Value2 -eq $oldValue) {
        $scratch.NumberFormat = "@"
        $scratch.Value = $newValue
        $scratch.Copy()
        $dst.PasteSpecial(-4163)   # xlPasteValues - value only, keep dst's style
        $scratch.Clear()           # leave the scratch cell exactly as found
    }
}
